$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying formatting (style) from the neighboring header cell G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add numeric values in H2:H3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
